$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This BOM sheet gets three new component rows (RESET button, R35-R38, and
# D014-D023), three previously-blank "Comment" cells get filled in for the
# connector rows (J5/J9, U2-U6, G1), and the R14/R15/... resistor group
# designator list is expanded. We reproduce this by inserting blank rows at
# the right spots (so surrounding rows keep their original cell styles) and
# then writing the new/changed values, using PasteSpecial(formats) to pull
# the correct look for the brand-new rows instead of the stray style that a
# plain Insert() right below the bold header row would otherwise create.
# ---------------------------------------------------------------------------

# 1) New "RESET" row, inserted as the new row 2 (pushes the rest down by one).
$ws.Rows("2:2").Insert()
$ws.Range("A3:D3").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

$ws.Range("A2").Value = "'431256083736"
$ws.Range("B2").Value = "RESET"
$ws.Range("C2").Value = "4312560837X6"

# Fill in the previously-empty Comment cells for the three connector rows
# that used to be rows 2-4 and are now rows 3-5.
$ws.Range("A3").Value = "6130XX21121_61300621121"
$ws.Range("A4").Value = "BTS7030-2EPA"
$ws.Range("A5").Value = "IFX27001TF_V50"

# 2) Expand the designator list for the 13K resistor row (old row 13, now row 14).
$ws.Range("B14").Value = "R14, R15, R16, R17, R19, R20, R23, R26, R31, R32, R33, R34, R41, R44, R53, R56"

# 3) New "1k2 / R35-R38" row, inserted as the new row 17.
$ws.Rows("17:17").Insert()
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)

$ws.Range("A17").Value = "1k2"
$ws.Range("B17").Value = "R35, R36, R37, R38"
$ws.Range("C17").Value = "R0805"

# 4) New "GS1J-L / D014-D023" row, inserted as the new row 31.
$ws.Rows("31:31").Insert()
$ws.Range("A30:C30").Copy()
$ws.Range("A31:C31").PasteSpecial(-4122)

$ws.Range("A31").Value = "GS1J-L"
$ws.Range("B31").Value = "D014, D016, D018, D019, D020, D021, D022, D023"
$ws.Range("C31").Value = "SMA_DO-214AC"

$ws.Range("C2:C33").Select()
